# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) marking Control (0) vs MDD (1) rows,
# and updates a handful of refit prediction/error values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in H1, matching the bold/bordered style used by the other headers.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Label column: 0 = Control, 1 = MDD, for both the 100- and 200-iteration blocks.
$labels = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0;  7  = 1;  8  = 1;  9  = 1;  10 = 1; 11 = 1;
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}
foreach ($row in $labels.Keys) {
    $ws.Range("H$row").Value = $labels[$row]
}

# Updated prediction / error values from the refit.
$ws.Range("D3").Value = 0.2766211362665506
$ws.Range("E3").Value = 0.2766211362665506

$ws.Range("D8").Value = 0.8432161281361807
$ws.Range("E8").Value = 0.1567838718638193

$ws.Range("D9").Value = 0.4742158276372366
$ws.Range("E9").Value = 0.5257841723627634

$ws.Range("F11").Value = 0.5883098244667053
